$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "예윤"
$ws.Range("B1").Value = "재현 서연 태훈"
$ws.Range("C1").Value = "예윤"
$ws.Range("D1").Value = "재현 예윤 혜지 태훈 한솔 서연"
$ws.Range("E1").Value = "재현 서연 태훈"
$ws.Range("A2").Value = "예윤 유진"
$ws.Range("B2").Value = "재현 태훈"
$ws.Range("C2").Value = "희지 유진"
$ws.Range("D2").Value = "병국 재현 예윤 현빈 희지 혜지 태훈 한솔 서연"
$ws.Range("E2").Value = "재현 태훈"
$ws.Range("A3").Value = "재현 혜지 유진"
$ws.Range("B3").Value = "태훈"
$ws.Range("C3").Value = "재현 혜지 유진"
$ws.Range("D3").Value = "한솔 서연"
$ws.Range("E3").Value = "병국 현빈 희지"
$ws.Range("A4").Value = "재현 혜지 서연"
$ws.Range("B4").Value = "태훈"
$ws.Range("C4").Value = "재현 혜지 서연"
$ws.Range("D4").Value = "준범 한솔 서연"
$ws.Range("E4").Value = "병국 예윤 현빈 희지 유진 한솔"
$ws.Range("A5").Value = "재현 혜지 서연 한솔"
$ws.Range("B5").Value = "예윤 유진 태훈"
$ws.Range("C5").Value = "재현 혜지 한솔 서연"
$ws.Range("D5").Value = "예윤 준범 유진 한솔 서연"
$ws.Range("E5").Value = "병국 예윤 현빈 희지 유진"
$ws.Range("A6").Value = "재현 희지 혜지 서연 한솔"
$ws.Range("B6").Value = "예윤 현빈 희지 혜지 유진 태훈"
$ws.Range("C6").Value = "병국 재현 희지 혜지 한솔 서연"
$ws.Range("D6").Value = "재현 예윤 현빈 희지 혜지 준범 유진 한솔 서연"
$ws.Range("E6").Value = "병국 예윤 현빈 희지 혜지 유진 태훈"
